$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select A9:K40 on the original sheet before copying, so the final saved
# selection/view state matches (A9 active, A9:K40 selected).
$ws1.Activate()
$ws1.Range("A9:K40").Select()
$excel.ActiveWindow.DisplayGridlines = $true

# Add a new worksheet right after "Data Harian - Table" and rename it.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "Sheet1"

# Copy the finalized daily-data table (header + 31 days) into the new
# sheet, starting at A1. This brings along values, shared-string refs and
# cell styles (header style + bordered/wrapped data style).
$src = $ws1.Range("A9:K40")
$src.Copy($newSheet.Range("A1"))

# The pasted data rows keep the "wrap text" cell style, which in the
# finalized workbook renders at a taller row height.
$newSheet.Rows("2:32").RowHeight = 28.8

# Selection on the new sheet covers the whole table.
$newSheet.Range("A1:K32").Select()

$newSheet.Activate()
